$wb = $excel.ActiveWorkbook

# --- Sheet "المدائح": append two new rows (57, 58) ---
$ws10 = $wb.Worksheets.Item("المدائح")
$ws10.Cells.Item(57,1).Value = "انظروا يديه تأملوا رجليه"
$ws10.Cells.Item(57,2).Value = "{E1A31474-7F00-4B7D-90FD-A67415E32872}"
$ws10.Cells.Item(57,3).Value = 1582
$ws10.Cells.Item(57,4).Value = 1596
$ws10.Cells.Item(57,5).Value = 15
$ws10.Cells.Item(58,1).Value = "يا ربنا مولى السلام"
$ws10.Cells.Item(58,2).Value = "{325B4DE6-FA63-4B6D-827F-BCC78EF1FEC9}"
$ws10.Cells.Item(58,3).Value = 1597
$ws10.Cells.Item(58,4).Value = 1639
$ws10.Cells.Item(58,5).Value = 43

# --- Sheet "التسبحة": insert a new entry before row 23, shifting subsequent rows,
#     which appends a brand-new row 67 at the end ---
$ws = $wb.Worksheets.Item("التسبحة")
$ws.Cells.Item(23,1).Value = 'ابصالية آدام للخمسين'
$ws.Cells.Item(23,2).Value = '{43AC03AD-AC75-480D-987F-66CB8CBE3883}'
$ws.Cells.Item(23,3).Value = 496
$ws.Cells.Item(23,4).Value = 519
$ws.Cells.Item(23,5).Value = 24
$ws.Cells.Item(24,1).Value = 'ابصالية آدام لعيد التجلي'
$ws.Cells.Item(24,2).Value = '{EF0F739B-A8DE-419D-8D45-757AA9347AB5}'
$ws.Cells.Item(24,3).Value = 520
$ws.Cells.Item(24,4).Value = 543
$ws.Cells.Item(24,5).Value = 24
$ws.Cells.Item(25,1).Value = 'ابصالية واطس لعيد التجلي'
$ws.Cells.Item(25,2).Value = '{95F02DE0-6540-4250-B6D4-213F4C9B73FC}'
$ws.Cells.Item(25,3).Value = 544
$ws.Cells.Item(25,4).Value = 567
$ws.Cells.Item(25,5).Value = 24
$ws.Cells.Item(26,1).Value = 'ابصالية آدام لصوم العذراء'
$ws.Cells.Item(26,2).Value = '{2908EF39-9CFE-4101-AED3-B54AD30D5A78}'
$ws.Cells.Item(26,3).Value = 568
$ws.Cells.Item(26,4).Value = 598
$ws.Cells.Item(26,5).Value = 31
$ws.Cells.Item(27,1).Value = 'ابصالية واطس لصوم العذراء'
$ws.Cells.Item(27,2).Value = '{222D1CFF-8162-4F43-A7FC-D6E04CE630E4}'
$ws.Cells.Item(27,3).Value = 599
$ws.Cells.Item(27,4).Value = 622
$ws.Cells.Item(27,5).Value = 24
$ws.Cells.Item(28,1).Value = 'ابصالية آدام لعيد العذراء'
$ws.Cells.Item(28,2).Value = '{CF62ACEE-48F9-4ABA-ADDC-6180BEC4873D}'
$ws.Cells.Item(28,3).Value = 623
$ws.Cells.Item(28,4).Value = 646
$ws.Cells.Item(28,5).Value = 24
$ws.Cells.Item(29,1).Value = 'ابصالية واطس لعيد العذراء'
$ws.Cells.Item(29,2).Value = '{E2D40FD7-171F-428B-86DB-65B332AB25F3}'
$ws.Cells.Item(29,3).Value = 647
$ws.Cells.Item(29,4).Value = 673
$ws.Cells.Item(29,5).Value = 27
$ws.Cells.Item(30,1).Value = 'تسبحة الأيام'
$ws.Cells.Item(30,2).Value = '{34B2EF10-1881-4BA7-95A7-7AB2F6F6651C}'
$ws.Cells.Item(30,3).Value = 674
$ws.Cells.Item(30,4).Value = 674
$ws.Cells.Item(30,5).Value = 1
$ws.Cells.Item(31,1).Value = 'ابصالية الأحد 1'
$ws.Cells.Item(31,2).Value = '{F8548FDB-8D40-484A-8D19-36EC50E838FD}'
$ws.Cells.Item(31,3).Value = 675
$ws.Cells.Item(31,4).Value = 698
$ws.Cells.Item(31,5).Value = 24
$ws.Cells.Item(32,1).Value = 'ابصالية الأحد الثانية'
$ws.Cells.Item(32,2).Value = '{F263153B-C7A8-4E6B-AACA-6F05AF050F2E}'
$ws.Cells.Item(32,3).Value = 699
$ws.Cells.Item(32,4).Value = 730
$ws.Cells.Item(32,5).Value = 32
$ws.Cells.Item(33,1).Value = 'ابصالية الاثنين'
$ws.Cells.Item(33,2).Value = '{C966C7AC-73F9-4177-AA7F-71D0428224AF}'
$ws.Cells.Item(33,3).Value = 731
$ws.Cells.Item(33,4).Value = 761
$ws.Cells.Item(33,5).Value = 31
$ws.Cells.Item(34,1).Value = 'ابصالية الثلاثاء'
$ws.Cells.Item(34,2).Value = '{31645468-E515-4F5A-85DB-DEE662F6432A}'
$ws.Cells.Item(34,3).Value = 762
$ws.Cells.Item(34,4).Value = 773
$ws.Cells.Item(34,5).Value = 12
$ws.Cells.Item(35,1).Value = 'ابصالية الأربعاء'
$ws.Cells.Item(35,2).Value = '{8ABA75EA-D793-46A0-8AE2-5B61A6B4FD7E}'
$ws.Cells.Item(35,3).Value = 774
$ws.Cells.Item(35,4).Value = 798
$ws.Cells.Item(35,5).Value = 25
$ws.Cells.Item(36,1).Value = 'ابصالية الخميس'
$ws.Cells.Item(36,2).Value = '{02352F94-02C4-4D7F-9247-697DA282E7C9}'
$ws.Cells.Item(36,3).Value = 799
$ws.Cells.Item(36,4).Value = 819
$ws.Cells.Item(36,5).Value = 21
$ws.Cells.Item(37,1).Value = 'ابصالية الجمعة'
$ws.Cells.Item(37,2).Value = '{E8504067-DC7B-4818-8157-B947A0A74D9A}'
$ws.Cells.Item(37,3).Value = 820
$ws.Cells.Item(37,4).Value = 839
$ws.Cells.Item(37,5).Value = 20
$ws.Cells.Item(38,1).Value = 'ابصالية السبت'
$ws.Cells.Item(38,2).Value = '{BF504610-6275-426C-A939-798A885C5C71}'
$ws.Cells.Item(38,3).Value = 840
$ws.Cells.Item(38,4).Value = 871
$ws.Cells.Item(38,5).Value = 32
$ws.Cells.Item(39,1).Value = 'مقدمة الثيؤطوكيات الأدام'
$ws.Cells.Item(39,2).Value = '{E358EDB7-F8FF-43DA-A8B6-81839E23E4C6}'
$ws.Cells.Item(39,3).Value = 872
$ws.Cells.Item(39,4).Value = 876
$ws.Cells.Item(39,5).Value = 5
$ws.Cells.Item(40,1).Value = 'مقدمة الثيؤطوكيات الواطس'
$ws.Cells.Item(40,2).Value = '{F96B080A-3FB2-430B-9BED-E692E913A9B0}'
$ws.Cells.Item(40,3).Value = 877
$ws.Cells.Item(40,4).Value = 878
$ws.Cells.Item(40,5).Value = 2
$ws.Cells.Item(41,1).Value = 'ثيؤطوكية الأحد 1-6'
$ws.Cells.Item(41,2).Value = '{8B50A9B8-162A-45FD-A40D-5405E501F1E6}'
$ws.Cells.Item(41,3).Value = 879
$ws.Cells.Item(41,4).Value = 963
$ws.Cells.Item(41,5).Value = 85
$ws.Cells.Item(42,1).Value = 'ثيؤطوكية الأحد 7-9'
$ws.Cells.Item(42,2).Value = '{64CE0420-479F-4F12-AE0C-F1218BF21635}'
$ws.Cells.Item(42,3).Value = 964
$ws.Cells.Item(42,4).Value = 1015
$ws.Cells.Item(42,5).Value = 52
$ws.Cells.Item(43,1).Value = 'ثيؤطوكية الأحد 10'
$ws.Cells.Item(43,2).Value = '{797B381C-F875-4BB4-8ACB-A5852FFBD8FC}'
$ws.Cells.Item(43,3).Value = 1016
$ws.Cells.Item(43,4).Value = 1021
$ws.Cells.Item(43,5).Value = 6
$ws.Cells.Item(44,1).Value = 'ثيؤطوكية الأحد 11-15'
$ws.Cells.Item(44,2).Value = '{67866127-A8D5-451C-B0C2-1CE6E6FBCD1F}'
$ws.Cells.Item(44,3).Value = 1022
$ws.Cells.Item(44,4).Value = 1070
$ws.Cells.Item(44,5).Value = 49
$ws.Cells.Item(45,1).Value = 'ثيؤطوكية الإثنين'
$ws.Cells.Item(45,2).Value = '{5022D768-2E12-4BEA-8D76-E3896BD58932}'
$ws.Cells.Item(45,3).Value = 1071
$ws.Cells.Item(45,4).Value = 1115
$ws.Cells.Item(45,5).Value = 45
$ws.Cells.Item(46,1).Value = 'ثيؤطوكية الثلاثاء'
$ws.Cells.Item(46,2).Value = '{D1BFEE47-99F3-4046-8C36-B6397205435B}'
$ws.Cells.Item(46,3).Value = 1116
$ws.Cells.Item(46,4).Value = 1163
$ws.Cells.Item(46,5).Value = 48
$ws.Cells.Item(47,1).Value = 'ثيؤطوكية الأربعاء'
$ws.Cells.Item(47,2).Value = '{5AD56D85-2906-43FB-98E9-FB96F1B37293}'
$ws.Cells.Item(47,3).Value = 1164
$ws.Cells.Item(47,4).Value = 1208
$ws.Cells.Item(47,5).Value = 45
$ws.Cells.Item(48,1).Value = 'ثيؤطوكية الخميس'
$ws.Cells.Item(48,2).Value = '{88249BFF-471A-47A1-B7BC-E5A5093EC8D7}'
$ws.Cells.Item(48,3).Value = 1209
$ws.Cells.Item(48,4).Value = 1304
$ws.Cells.Item(48,5).Value = 96
$ws.Cells.Item(49,1).Value = 'ثيؤطوكية الجمعة'
$ws.Cells.Item(49,2).Value = '{6C9361D4-74F3-4201-B28D-7EB59C9D9A46}'
$ws.Cells.Item(49,3).Value = 1305
$ws.Cells.Item(49,4).Value = 1333
$ws.Cells.Item(49,5).Value = 29
$ws.Cells.Item(50,1).Value = 'ثيؤطوكية السبت'
$ws.Cells.Item(50,2).Value = '{25CBC7C4-A68C-4EBD-B127-98DA707B3413}'
$ws.Cells.Item(50,3).Value = 1334
$ws.Cells.Item(50,4).Value = 1370
$ws.Cells.Item(50,5).Value = 37
$ws.Cells.Item(51,1).Value = 'ثيؤطوكية الأحد 16-18'
$ws.Cells.Item(51,2).Value = '{4E843BF3-1D30-4BED-905C-E66AA3D90EC5}'
$ws.Cells.Item(51,3).Value = 1371
$ws.Cells.Item(51,4).Value = 1383
$ws.Cells.Item(51,5).Value = 13
$ws.Cells.Item(52,1).Value = 'لبش الإثنين'
$ws.Cells.Item(52,2).Value = '{B08DAA27-DC93-470F-8EE4-DBA2CDED73FF}'
$ws.Cells.Item(52,3).Value = 1384
$ws.Cells.Item(52,4).Value = 1395
$ws.Cells.Item(52,5).Value = 12
$ws.Cells.Item(53,1).Value = 'لبش الثلاثاء'
$ws.Cells.Item(53,2).Value = '{EA64C1D5-8011-4ED1-AECA-ACA0D1D96925}'
$ws.Cells.Item(53,3).Value = 1396
$ws.Cells.Item(53,4).Value = 1404
$ws.Cells.Item(53,5).Value = 9
$ws.Cells.Item(54,1).Value = 'لبش الأربعاء'
$ws.Cells.Item(54,2).Value = '{824D594F-C079-4552-882A-CC297F319D7D}'
$ws.Cells.Item(54,3).Value = 1405
$ws.Cells.Item(54,4).Value = 1419
$ws.Cells.Item(54,5).Value = 15
$ws.Cells.Item(55,1).Value = 'لبش الخميس'
$ws.Cells.Item(55,2).Value = '{260E1FAC-A9F6-4E94-BAAB-EFD045CD242D}'
$ws.Cells.Item(55,3).Value = 1420
$ws.Cells.Item(55,4).Value = 1436
$ws.Cells.Item(55,5).Value = 17
$ws.Cells.Item(56,1).Value = 'لبش الجمعة'
$ws.Cells.Item(56,2).Value = '{27A6E4EC-9C9A-4029-8EAF-A984FA647997}'
$ws.Cells.Item(56,3).Value = 1437
$ws.Cells.Item(56,4).Value = 1456
$ws.Cells.Item(56,5).Value = 20
$ws.Cells.Item(57,1).Value = 'شيرات السبت 1'
$ws.Cells.Item(57,2).Value = '{FA5AF629-FC64-4123-92EA-193DFE2229CC}'
$ws.Cells.Item(57,3).Value = 1457
$ws.Cells.Item(57,4).Value = 1468
$ws.Cells.Item(57,5).Value = 12
$ws.Cells.Item(58,1).Value = 'شيرات السبت 2'
$ws.Cells.Item(58,2).Value = '{2DF7B6FE-B056-4813-B72C-DFE470371815}'
$ws.Cells.Item(58,3).Value = 1469
$ws.Cells.Item(58,4).Value = 1484
$ws.Cells.Item(58,5).Value = 16
$ws.Cells.Item(59,1).Value = 'مقدمة الدفنار'
$ws.Cells.Item(59,2).Value = '{DBBEB49F-3396-41D0-81FF-0A028C3CB4DA}'
$ws.Cells.Item(59,3).Value = 1485
$ws.Cells.Item(59,4).Value = 1485
$ws.Cells.Item(59,5).Value = 1
$ws.Cells.Item(60,1).Value = 'مقدمة الدفنار الآدام'
$ws.Cells.Item(60,2).Value = '{0420AA0C-B21A-478D-88EA-8378E9539EDE}'
$ws.Cells.Item(60,3).Value = 1486
$ws.Cells.Item(60,4).Value = 1489
$ws.Cells.Item(60,5).Value = 4
$ws.Cells.Item(61,1).Value = 'مقدمة الدفنار الواطس'
$ws.Cells.Item(61,2).Value = '{F2541D73-C210-4196-BE50-DF6E6142A86C}'
$ws.Cells.Item(61,3).Value = 1490
$ws.Cells.Item(61,4).Value = 1493
$ws.Cells.Item(61,5).Value = 4
$ws.Cells.Item(62,1).Value = 'الدفنار'
$ws.Cells.Item(62,2).Value = '{A509B738-02BB-455A-944E-9E56D85C8942}'
$ws.Cells.Item(62,3).Value = 1494
$ws.Cells.Item(62,4).Value = 1495
$ws.Cells.Item(62,5).Value = 2
$ws.Cells.Item(63,1).Value = 'ختام الثؤطوكيات الادام'
$ws.Cells.Item(63,2).Value = '{14A3A43C-A9F7-45A8-A510-EE3F33D99572}'
$ws.Cells.Item(63,3).Value = 1496
$ws.Cells.Item(63,4).Value = 1512
$ws.Cells.Item(63,5).Value = 17
$ws.Cells.Item(64,1).Value = 'ختام الثيؤطوكيات الواطس'
$ws.Cells.Item(64,2).Value = '{BF439D71-64D1-4376-8E4A-812437425EBB}'
$ws.Cells.Item(64,3).Value = 1513
$ws.Cells.Item(64,4).Value = 1532
$ws.Cells.Item(64,5).Value = 20
$ws.Cells.Item(65,1).Value = 'قانون الايمان'
$ws.Cells.Item(65,2).Value = '{A12368B5-4E89-4682-AF79-DC1979BA120B}'
$ws.Cells.Item(65,3).Value = 1533
$ws.Cells.Item(65,4).Value = 1547
$ws.Cells.Item(65,5).Value = 15
$ws.Cells.Item(66,1).Value = 'ختام التسبحة'
$ws.Cells.Item(66,2).Value = '{BABBF91F-DC4B-4DBC-A6C1-054AEA7290F3}'
$ws.Cells.Item(66,3).Value = 1548
$ws.Cells.Item(66,4).Value = 1562
$ws.Cells.Item(66,5).Value = 15
$ws.Cells.Item(67,1).Value = 'قدوس قدوس قدوس'
$ws.Cells.Item(67,2).Value = '{F2F363F3-5DD8-474B-94A0-6895758AB76D}'
$ws.Cells.Item(67,3).Value = 1563
$ws.Cells.Item(67,4).Value = 1567
$ws.Cells.Item(67,5).Value = 5
